# UI Minor update and Bug Fix List Update
# Adds four new bug-tracker rows (PIUTANG / HUTANG payment modules) to the
# "11 JAN 2017" sheet, then refreshes column B's width so the longer new
# module names fit, and leaves the selection on the newly added blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New bug rows -------------------------------------------------------
# Write order matches how the shared-string table grows (new unique
# strings are appended in first-use order), so it mirrors the original
# authoring sequence: both PIUTANG rows' labels/description first, then
# back-fill the CUSTOMER row's recurring error text, then the HUTANG rows.
$ws.Range("B35").Value = "PEMBAYARAN PIUTANG PER-CUSTOMER"
$ws.Range("B37").Value = "PEMBAYARAN PIUTANG PER-INVOICE"
$ws.Range("C37").Value = "Pada saat save error mysqlexception unknown column ORIGIN_SALES_INVOICE"
$ws.Range("C35").Value = "error nullreferenceexception pada prosedur unregisterglobalhotkey"

$ws.Range("B39").Value = "PEMBAYARAN HUTANG PER-INVOICE"
$ws.Range("C39").Value = "error nullreferenceexception pada prosedur unregisterglobalhotkey"

$ws.Range("B41").Value = "PEMBAYARAN HUTANG PER-CUSTOMER"
$ws.Range("C41").Value = "error nullreferenceexception pada prosedur unregisterglobalhotkey"

# --- Column B is now widest for "PEMBAYARAN PIUTANG PER-CUSTOMER"; -----
# widen/best-fit it like Excel would after typing the longer labels.
$ws.Columns("B:B").ColumnWidth = 36

# --- Leave the selection on the next blank row, matching where the ----
# user was working after the last entry.
[void]$ws.Range("A40:XFD40").Select()
